$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Replace the row-11 rule label ("R40") with the text "1".
# A leading apostrophe forces Excel to store it as text (matching the
# original cell, which also held a text/string value) rather than
# re-interpreting "1" as a number.
$ws.Range("B11").Value = "'1"
